$wb = $excel.ActiveWorkbook

# --- "Boolean" sheet: expand the two combined "QaZ" rows into per-mode files ---
$wsBoolean = $wb.Worksheets.Item("Boolean")

# Row 17 held "trans/BVTQaZ/BVTQaZ.csv" -> split into 6 mode-specific files
$wsBoolean.Rows("18:22").Insert()
$wsBoolean.Range("A17").Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$wsBoolean.Range("A18").Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$wsBoolean.Range("A19").Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$wsBoolean.Range("A20").Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$wsBoolean.Range("A21").Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$wsBoolean.Range("A22").Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# After the insert above, "trans/VTQaZ/VTQaZ.csv" (originally row 21) is now row 26.
# Split it into 6 mode-specific files as well.
$wsBoolean.Rows("27:31").Insert()
$wsBoolean.Range("A26").Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
$wsBoolean.Range("A27").Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$wsBoolean.Range("A28").Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$wsBoolean.Range("A29").Value = "trans/VTQaZ/VTQaZ-rail.csv"
$wsBoolean.Range("A30").Value = "trans/VTQaZ/VTQaZ-ships.csv"
$wsBoolean.Range("A31").Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# --- View/selection state ---
# "Integer" sheet: remember a selection at A13 (no longer the active tab)
$wsInteger = $wb.Worksheets.Item("Integer")
$wsInteger.Activate() | Out-Null
$wsInteger.Range("A13").Select() | Out-Null

# "Boolean" sheet: remember a selection at A32 (scrolled down to the new rows)
$wsBoolean.Activate() | Out-Null
$wsBoolean.Range("A32").Select() | Out-Null

# "About" sheet becomes the active tab/selected sheet on save
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate() | Out-Null
$wsAbout.Range("A1").Select() | Out-Null
